$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Browser Name column (B2:B7) from "Chrome" to "Firefox"
$ws.Range("B2:B7").Value = "Firefox"

# Update the active cell / selection on the sheet view to B12
$ws.Activate()
$ws.Range("B12").Select()
